# Regenerate merged AHB files
# ----------------------------------------------------------------------------
# This script reproduces (via Excel COM automation) the change described by
# the target diff:
#   1. Rename the shared header strings:
#        "<Name>_old" -> "<Name>_FV2404"   (columns A..J)
#        "<Name>_new" -> "<Name>_FV2410"   (columns L..U)
#      (column K, "diff", is left untouched)
#   2. Turn the used range A1:U81 into a native Excel Table ("Table1") with
#      an AutoFilter, 21 columns named after the (renamed) header row.
#   3. Freeze the header row (pane split under row 1).
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header labels in place -----------------------------------
# Using Find/Replace (substring match) on the whole sheet touches only the
# header cells in row 1, since "_old"/"_new" do not occur anywhere else in
# the data.
$ws.Cells.Replace("_old", "_FV2404")
$ws.Cells.Replace("_new", "_FV2410")

# --- 2. Build the table ------------------------------------------------------
# Excel would normally bake the header row's existing formatting into a new
# "dxf" (and reference it via headerRowDxfId) when a table is created on a
# range whose header already carries explicit cell formatting. The target
# workbook does not have that extra dxf, so we temporarily neutralise the
# header formatting (saving it first onto a scratch row far outside the used
# range) before creating the table, then restore the original formatting
# (reusing the very same style, so no new style record is created) and wipe
# the scratch row again.
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A200:U200")

$headerRange.Copy()
$scratch.PasteSpecial(-4122)            # xlPasteFormats
$excel.CutCopyMode = 0

$headerRange.Style = "Normal"

$tableRange = $ws.Range("A1:U81")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)   # xlSrcRange, xlYes
$tbl.TableStyle = ""                     # no explicit table style, like the target

$scratch.Copy()
$headerRange.PasteSpecial(-4122)        # xlPasteFormats
$excel.CutCopyMode = 0

$scratch.Clear()

# --- 3. Freeze the header row -------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Applied header rename, table creation, and freeze pane."
